$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 825.2192
$ws.Range("I15").Value = 825.2192
$ws.Range("K15").Value = 2475.6576
$ws.Range("M15").Value = -2306.6576
$ws.Range("H17").Value = 890730.06
$ws.Range("J17").Value = 890730.06
$ws.Range("L17").Value = 2672190.18
$ws.Range("N17").Value = -2672526.18
$ws.Range("H113").Value = 3962.0
$ws.Range("J113").Value = 6609.4
$ws.Range("L113").Value = 6609.4
$ws.Range("N113").Value = -13117.4
$ws.Range("H138").Value = 10105855.0
$ws.Range("J138").Value = 17551260.0
$ws.Range("L138").Value = 52653780.0
$ws.Range("N138").Value = -52664060.0
$ws.Range("H141").Value = 1913.762
$ws.Range("I141").Value = 1431.0
$ws.Range("K141").Value = 4293.0
$ws.Range("M141").Value = 887.0

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 742.6667
$ws.Range("J2").Value = 2006.5
$ws.Range("L2").Value = 2006.5
$ws.Range("N2").Value = -2232.5
$ws.Range("H32").Value = 17550190.0
$ws.Range("I32").Value = 21280112.0
$ws.Range("J32").Value = 19559.9
$ws.Range("K32").Value = 21280112.0
$ws.Range("L32").Value = 19559.9
$ws.Range("M32").Value = -21279825.0
$ws.Range("N32").Value = -20133.9
$ws.Range("H45").Value = 2886.125
$ws.Range("I45").Value = 2012.0
$ws.Range("K45").Value = 2012.0
$ws.Range("M45").Value = -1635.0
$ws.Range("H116").Value = 742.6667
$ws.Range("J116").Value = 2006.5
$ws.Range("L116").Value = 2006.5
$ws.Range("N116").Value = -6594.5
$ws.Range("H122").Value = 3008.7144
$ws.Range("I122").Value = 1956.5555
$ws.Range("K122").Value = 5869.666499999999
$ws.Range("M122").Value = -3419.666499999999
$ws.Range("H132").Value = 111115496.0
$ws.Range("I132").Value = 4932.0
$ws.Range("K132").Value = 14796.0
$ws.Range("M132").Value = -12266.0
$ws.Range("H135").Value = 45353.0
$ws.Range("J135").Value = 45353.0
$ws.Range("L135").Value = 45353.0
$ws.Range("N135").Value = -55493.0

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 742.6667
$ws.Range("J3").Value = 2006.5
$ws.Range("L3").Value = 2006.5
$ws.Range("N3").Value = -2234.5
$ws.Range("H86").Value = 12721.321
$ws.Range("I86").Value = 11458.7
$ws.Range("J86").Value = 13422.777
$ws.Range("K86").Value = 11458.7
$ws.Range("L86").Value = 13422.777
$ws.Range("M86").Value = -10335.7
$ws.Range("N86").Value = -15668.777
$ws.Range("H89").Value = 12721.321
$ws.Range("I89").Value = 11458.7
$ws.Range("J89").Value = 13422.777
$ws.Range("K89").Value = 57293.5
$ws.Range("L89").Value = 67113.885
$ws.Range("M89").Value = -51677.5
$ws.Range("N89").Value = -78345.885
$ws.Range("H134").Value = 3200.973
$ws.Range("J134").Value = 11499.0
$ws.Range("L134").Value = 34497.0
$ws.Range("N134").Value = -39567.0

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 674.3
$ws.Range("I16").Value = 488.66666
$ws.Range("K16").Value = 488.66666
$ws.Range("M16").Value = -201.66666
$ws.Range("H31").Value = 21280762.0
$ws.Range("I31").Value = 3148.5625
$ws.Range("J31").Value = 66673004.0
$ws.Range("K31").Value = 3148.5625
$ws.Range("L31").Value = 66673004.0
$ws.Range("M31").Value = -2853.5625
$ws.Range("N31").Value = -66673594.0
$ws.Range("H34").Value = 21280762.0
$ws.Range("I34").Value = 3148.5625
$ws.Range("J34").Value = 66673004.0
$ws.Range("K34").Value = 3148.5625
$ws.Range("L34").Value = 66673004.0
$ws.Range("M34").Value = -2946.5625
$ws.Range("N34").Value = -66673408.0
$ws.Range("H99").Value = 12056.629
$ws.Range("I99").Value = 13180.182
$ws.Range("J99").Value = 11541.667
$ws.Range("K99").Value = 13180.182
$ws.Range("L99").Value = 11541.667
$ws.Range("M99").Value = -11682.182
$ws.Range("N99").Value = -14537.667
$ws.Range("H113").Value = 674.3
$ws.Range("I113").Value = 488.66666
$ws.Range("K113").Value = 488.66666
$ws.Range("M113").Value = 1681.33334
$ws.Range("H126").Value = 12056.629
$ws.Range("I126").Value = 13180.182
$ws.Range("J126").Value = 11541.667
$ws.Range("K126").Value = 39540.546
$ws.Range("L126").Value = 34625.001
$ws.Range("M126").Value = -37070.546
$ws.Range("N126").Value = -39565.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H116").Value = 1749.6666
$ws.Range("I116").Value = 499.33334
$ws.Range("K116").Value = 1498.00002
$ws.Range("M116").Value = 1943.99998
$ws.Range("H131").Value = 2794.1785
$ws.Range("I131").Value = 1658.6666
$ws.Range("J131").Value = 3103.8635
$ws.Range("K131").Value = 4975.9998
$ws.Range("L131").Value = 9311.5905
$ws.Range("M131").Value = 64.0002000000004
$ws.Range("N131").Value = -19391.5905
$ws.Range("H133").Value = 11120.728
$ws.Range("J133").Value = 14850.923
$ws.Range("L133").Value = 44552.769
$ws.Range("N133").Value = -54672.769
$ws.Range("H134").Value = 3729.9443
$ws.Range("I134").Value = 1604.9333
$ws.Range("K134").Value = 4814.7999
$ws.Range("M134").Value = 255.2001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 752.94116
$ws.Range("I107").Value = 656.0
$ws.Range("J107").Value = 862.0
$ws.Range("K107").Value = 656.0
$ws.Range("L107").Value = 862.0
$ws.Range("M107").Value = 1264.0
$ws.Range("N107").Value = -4702.0
$ws.Range("H122").Value = 5057.826
$ws.Range("I122").Value = 2031.091
$ws.Range("J122").Value = 7832.3335
$ws.Range("K122").Value = 6093.272999999999
$ws.Range("L122").Value = 23497.0005
$ws.Range("M122").Value = -3643.272999999999
$ws.Range("N122").Value = -28397.0005
$ws.Range("H126").Value = 30007794.0
$ws.Range("I126").Value = 33346324.0
$ws.Range("J126").Value = 28576996.0
$ws.Range("K126").Value = 100038972.0
$ws.Range("L126").Value = 85730988.0
$ws.Range("M126").Value = -100036502.0
$ws.Range("N126").Value = -85735928.0

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 2235.8276
$ws.Range("J22").Value = 2657.5264
$ws.Range("L22").Value = 2657.5264
$ws.Range("N22").Value = -3247.5264
$ws.Range("H27").Value = 2235.8276
$ws.Range("J27").Value = 2657.5264
$ws.Range("L27").Value = 2657.5264
$ws.Range("N27").Value = -2871.5264
$ws.Range("H40").Value = 3793.7273
$ws.Range("I40").Value = 3973.1
$ws.Range("J40").Value = 2000.0
$ws.Range("K40").Value = 3973.1
$ws.Range("L40").Value = 2000.0
$ws.Range("M40").Value = -3837.1
$ws.Range("N40").Value = -2272.0
$ws.Range("H46").Value = 1185.8636
$ws.Range("I46").Value = 670.35297
$ws.Range("J46").Value = 2938.6
$ws.Range("K46").Value = 670.35297
$ws.Range("L46").Value = 2938.6
$ws.Range("M46").Value = -482.35297
$ws.Range("N46").Value = -3314.6
$ws.Range("H82").Value = 7398.8
$ws.Range("I82").Value = 5000.0
$ws.Range("J82").Value = 7998.5
$ws.Range("K82").Value = 5000.0
$ws.Range("L82").Value = 7998.5
$ws.Range("M82").Value = -4639.0
$ws.Range("N82").Value = -8720.5
$ws.Range("H85").Value = 7398.8
$ws.Range("I85").Value = 5000.0
$ws.Range("J85").Value = 7998.5
$ws.Range("K85").Value = 5000.0
$ws.Range("L85").Value = 7998.5
$ws.Range("M85").Value = -3752.0
$ws.Range("N85").Value = -10494.5
$ws.Range("H100").Value = 3931.111
$ws.Range("I100").Value = 2997.0
$ws.Range("K100").Value = 2997.0
$ws.Range("M100").Value = -2456.0
$ws.Range("H106").Value = 26634.5
$ws.Range("J106").Value = 26634.5
$ws.Range("L106").Value = 26634.5
$ws.Range("N106").Value = -29158.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H3").Value = 15000.0
$ws.Range("I3").Value = 15000.0
$ws.Range("K3").Value = 15000.0
$ws.Range("M3").Value = -14886.0
$ws.Range("H113").Value = 715.7727
$ws.Range("I113").Value = 312.93332
$ws.Range("K113").Value = 938.7999599999999
$ws.Range("M113").Value = 1231.20004
$ws.Range("H122").Value = 43480028.0
$ws.Range("I122").Value = 50001480.0
$ws.Range("J122").Value = 3684.6667
$ws.Range("K122").Value = 150004440.0
$ws.Range("L122").Value = 11054.0001
$ws.Range("M122").Value = -150001990.0
$ws.Range("N122").Value = -15954.0001
$ws.Range("H126").Value = 3362.4546
$ws.Range("I126").Value = 4271.0835
$ws.Range("J126").Value = 939.44446
$ws.Range("K126").Value = 12813.2505
$ws.Range("L126").Value = 2818.33338
$ws.Range("M126").Value = -10343.2505
$ws.Range("N126").Value = -7758.33338
